$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells D1/E1: Param2/Param3 -> foo/bar
$ws.Range("D1").Value = "foo"
$ws.Range("E1").Value = "bar"

# Update numeric values in columns D and E for rows 2-5
$ws.Range("D2").Value = 0.2
$ws.Range("E2").Value = 0.1

$ws.Range("D3").Value = 0.1
$ws.Range("E3").Value = 0.9

$ws.Range("D4").Value = 0.9
$ws.Range("E4").Value = 0.9

$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 0.2

# Update selection to D1:E5 with active cell D1
$ws.Range("D1:E5").Select()
